# "Generate Report for Archive"
# The localization status report is regenerated: the shared "Status" value
# used across the Overview / zh-cn / de-de sheets moves from
# "Ready for handoff" to "In Translation", and the "Status" columns are
# re-sized (narrower) to fit the refreshed value.

$wb = $excel.ActiveWorkbook
$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

# --- zh-cn sheet: Status column (C) --------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de sheet: Status column (C) --------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus

# --- Resize the "Status" columns to match the refreshed report layout ----
# (engine quantizes ColumnWidth to the nearest 1/6 character; 12.5 is the
# closest achievable value to the target width.)
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn" status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de" status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
